# Updates the cryptocurrency price/volume table on Sheet1 (rows 2-51)
# to reflect the latest scrape, matching the commit:
# "Updated cryptos list on Tue May 21 13:23:46 UTC 2024 with GitHub Actions"
#
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
# Column D holds price text that often looks like a plain number
# (e.g. "1.00", "0.542") even though the sheet stores it as text, so a
# leading apostrophe is used where needed to force Excel to keep it as
# text instead of auto-converting it to a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.036.02'
$ws.Range("E2").Value = '  +5.92%  '
$ws.Range("D3").Value = '3.781.62'
$ws.Range("E3").Value = '  +22.41%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''611.35'
$ws.Range("E5").Value = '  +6.92%  '
$ws.Range("D6").Value = '''178.72'
$ws.Range("E6").Value = '  +0.95%  '
$ws.Range("D7").Value = '3.780.28'
$ws.Range("E7").Value = '  +22.38%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '''0.542'
$ws.Range("E9").Value = '  +5.83%  '
$ws.Range("D10").Value = '''0.167'
$ws.Range("E10").Value = '  +10.16%  '
$ws.Range("D11").Value = '''6.40'
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("D12").Value = '''0.498'
$ws.Range("E12").Value = '  +7.04%  '
$ws.Range("D13").Value = '''40.62'
$ws.Range("E13").Value = '  +13.38%  '
$ws.Range("D14").Value = '''0.0000257'
$ws.Range("E14").Value = '  +7.01%  '
$ws.Range("D15").Value = '4.412.49'
$ws.Range("E15").Value = '  +22.41%  '
$ws.Range("D16").Value = '3.786.16'
$ws.Range("E16").Value = '  +22.57%  '
$ws.Range("D17").Value = '71.247.41'
$ws.Range("E17").Value = '  +6.35%  '
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("D19").Value = '''7.52'
$ws.Range("E19").Value = '  +7.67%  '
$ws.Range("D20").Value = '''523.72'
$ws.Range("E20").Value = '  +7.50%  '
$ws.Range("D21").Value = '''16.72'
$ws.Range("E21").Value = '  +1.42%  '
$ws.Range("D22").Value = '''9.43'
$ws.Range("E22").Value = '  +23.14%  '
$ws.Range("D23").Value = '''0.743'
$ws.Range("E23").Value = '  +8.72%  '
$ws.Range("D24").Value = '''88.36'
$ws.Range("E24").Value = '  +6.08%  '
$ws.Range("D25").Value = '''2.48'
$ws.Range("E25").Value = '  +10.26%  '
$ws.Range("D26").Value = '''13.48'
$ws.Range("E26").Value = '  +7.12%  '
$ws.Range("D27").Value = '''11.00'
$ws.Range("E27").Value = '  +8.17%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").Value = '  +28.62%  '
$ws.Range("D30").Value = '''2.51'
$ws.Range("E30").Value = '  +10.06%  '
$ws.Range("D31").Value = '''2.91'
$ws.Range("E31").Value = '  +12.93%  '
$ws.Range("D32").Value = '''8.00'
$ws.Range("E32").Value = '  +2.01%  '
$ws.Range("D33").Value = '''32.20'
$ws.Range("E33").Value = '  +15.17%  '
$ws.Range("E34").Value = '  +3.54%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("D36").Value = '''1.06'
$ws.Range("E36").Value = '  +11.93%  '
$ws.Range("D37").Value = '''6.11'
$ws.Range("E37").Value = '  +9.95%  '
$ws.Range("D38").Value = '''2.22'
$ws.Range("E38").Value = '  +10.85%  '
$ws.Range("D39").Value = '''0.339'
$ws.Range("E39").Value = '  +9.24%  '
$ws.Range("E40").Value = '  +6.93%  '
$ws.Range("D41").Value = '''51.59'
$ws.Range("E41").Value = '  +5.24%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '3.152.30'
$ws.Range("E42").Value = '  +13.16%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = '''427.52'
$ws.Range("E43").Value = '  +16.40%  '
$ws.Range("D44").Value = '''8.80'
$ws.Range("E44").Value = '  +7.33%  '
$ws.Range("D45").Value = '''44.05'
$ws.Range("E45").Value = '  -7.04%  '
$ws.Range("D46").Value = '''2.78'
$ws.Range("E46").Value = '  +2.02%  '
$ws.Range("E47").Value = '  +6.98%  '
$ws.Range("D48").Value = '''27.76'
$ws.Range("E48").Value = '  +9.17%  '
$ws.Range("D49").Value = '''141.36'
$ws.Range("E49").Value = '  +4.77%  '
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").Value = '''2.47'
$ws.Range("E51").Value = '  +7.07%  '

